$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dict_sheet")

# The key/value pairs were scrambled (rows out of order relative to each
# other). Fix the data so each row holds the matching key/value pair in
# the correct order: key1/value1, key2/value2, key3/value3.
$ws.Range("A1").Value = "key1"
$ws.Range("B1").Value = "value1"
$ws.Range("A2").Value = "key2"
$ws.Range("B2").Value = "value2"
$ws.Range("A3").Value = "key3"
$ws.Range("B3").Value = "value3"
